$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) --------------------------------------------
$ws.Range("AA1").Value = "date_data_collection"
$ws.Range("AB1").Value = "identification"

# --- New data cells --------------------------------------------------------
# AA2 holds the text "0505" (looks numeric, so it must stay text). Typing it
# directly makes Excel auto-detect a number; instead write it with a leading
# apostrophe, then copy the number format/style from a plain left-aligned
# text cell (A2) on top of it so the final style matches the rest of the
# column instead of getting a fresh "quote prefix" style.
$ws.Range("AA2").Value = "'0505"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("AA2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AB2").Value = 1
$ws.Range("AB3").Value = 1
$ws.Range("AB4").Value = 1
$ws.Range("AB5").Value = 1

# --- Column width for the new column AA (27) --------------------------------
# (COM ColumnWidth and the stored OOXML <col width> differ by the sheet's
# "Maximum Digit Width" padding; 22.95 is the closest COM value that rounds
# back to the ~23.83-character stored width used by the rest of the sheet.)
$ws.Columns.Item(27).ColumnWidth = 22.95

# --- Selection / view state --------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("AB6").Select()
